{"js": "// Change \"All responses will be done on Blackboard. \" to\n// \"All responses will be done on reddit, see: \" followed by a\n// HYPERLINK field (displaying https://www.reddit.com/r/QCMath2903/)\n// and the trailing text \" after I add the class.\"\n\n// Step 1: swap \"Blackboard. \" for \"reddit, see: \" in place so the\n// surrounding run keeps its original character formatting (rFonts,\n// sz, rtl, lang, ...). insertText with InsertLocation.replace edits\n// the text node of the existing run rather than fabricating a new one.\nconst target = context.document.body.search(\"Blackboard. \", { matchCase: true, matchWholeWord: false });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error('Could not find \"Blackboard. \" in the document body.');\n}\n\nconst blackboardRange = target.items[0];\nblackboardRange.insertText(\"reddit, see: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: locate the (now unique) text we just inserted and append the\n// hyperlink field + closing sentence right after it, via a flat-OPC\n// OOXML fragment so we get the exact begin/instrText/separate/text/end\n// field-code run sequence (with the Hyperlink.0 character style on the\n// field runs) rather than Office.js's own higher-level hyperlink API.\nconst anchor = context.document.body.search(\"reddit, see: \", { matchCase: true, matchWholeWord: false });\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error('Could not find \"reddit, see: \" after the text replacement.');\n}\n\nconst anchorRange = anchor.items[0];\n\nconst hyperlinkUrl = \"https://www.reddit.com/r/QCMath2903/\";\n\nconst flatOpcFragment =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r>' +\n  '<w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:fldChar w:fldCharType=\"begin\" w:fldLock=\"0\"/>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:instrText xml:space=\"preserve\"> HYPERLINK \"' + hyperlinkUrl + '\"</w:instrText>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:fldChar w:fldCharType=\"separate\" w:fldLock=\"0\"/>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:rtl w:val=\"0\"/><w:lang w:val=\"en-US\"/></w:rPr>' +\n  '<w:t>' + hyperlinkUrl + '</w:t>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>' +\n  '<w:fldChar w:fldCharType=\"end\" w:fldLock=\"0\"/>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:rtl w:val=\"0\"/><w:lang w:val=\"en-US\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\"> after I add the class.</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nanchorRange.insertOoxml(flatOpcFragment, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Change \"All responses will be done on Blackboard. \" to\n# \"All responses will be done on reddit, see: \" followed by a HYPERLINK\n# field (displaying https://www.reddit.com/r/QCMath2903/) and the\n# trailing sentence \" after I add the class.\"\n\n$d = $word.ActiveDocument\n\n# Step 1: swap \"Blackboard. \" for \"reddit, see: \" in place so the\n# surrounding run keeps its original character formatting (rFonts, sz,\n# rtl, lang, ...). Assigning Range.Text edits the existing run's text\n# rather than fabricating a new one.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Blackboard. \")\nif (-not $found) {\n    throw 'Could not find \"Blackboard. \" in the document body.'\n}\n$findRange.Text = \"reddit, see: \"\n\n# Step 2: re-locate the (now unique) text we just inserted and collapse\n# to its end, then splice in the hyperlink field + closing sentence via\n# a flat-OPC OOXML fragment so we get the exact\n# begin/instrText/separate/text/end field-code run sequence (with the\n# Hyperlink.0 character style on the field runs), matching what Word\n# itself writes for Insert > Link.\n$insertRange = $d.Content\n$found2 = $insertRange.Find.Execute(\"reddit, see: \")\nif (-not $found2) {\n    throw 'Could not find \"reddit, see: \" after the text replacement.'\n}\n$insertRange.SetRange($insertRange.End, $insertRange.End)\n\n$hyperlinkUrl = \"https://www.reddit.com/r/QCMath2903/\"\n\n$xml = '<?xml version=\"1.0\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n'<pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n'<w:body>' + `\n'<w:p>' + `\n'<w:r><w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:fldChar w:fldCharType=\"begin\" w:fldLock=\"0\"/></w:r>' + `\n('<w:r><w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:instrText xml:space=\"preserve\"> HYPERLINK \"' + $hyperlinkUrl + '\"</w:instrText></w:r>') + `\n'<w:r><w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:fldChar w:fldCharType=\"separate\" w:fldLock=\"0\"/></w:r>' + `\n('<w:r><w:rPr><w:rStyle w:val=\"Hyperlink.0\"/><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:rtl w:val=\"0\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>' + $hyperlinkUrl + '</w:t></w:r>') + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:fldChar w:fldCharType=\"end\" w:fldLock=\"0\"/></w:r>' + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Calibri\" w:cs=\"Calibri\" w:hAnsi=\"Calibri\" w:eastAsia=\"Calibri\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:rtl w:val=\"0\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> after I add the class.</w:t></w:r>' + `\n'</w:p>' + `\n'</w:body>' + `\n'</w:document>' + `\n'</pkg:xmlData>' + `\n'</pkg:part>' + `\n'</pkg:package>'\n\n$insertRange.InsertXML($xml, \"End\")\n"}
